$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header: "Table ID" -> "Table_ID [Anguss_tool]" ---
$ws.Range("B1").Value = "Table_ID [Anguss_tool]"

# --- 2. Rebuild the "Bearbeiten Außenseite..." instructional text with a
#        partially-bold run around [Anguss_note], and make C5 match the
#        same corrected text (it previously held a mangled variant). ---
$bearbeitenText = "Bearbeiten Außenseite:- Anguss Ansatz [Anguss_note] entfernen."
$boldStart = 40
$boldLen = 11

$c2 = $ws.Range("C2")
$c2.Value = $bearbeitenText

$pre = $c2.Characters(1, ($boldStart - 1))
$pre.Font.Name = "Open Sans"
$pre.Font.Size = 10
$pre.Font.Bold = $false

$mid = $c2.Characters($boldStart, $boldLen)
$mid.Font.Name = "Open Sans"
$mid.Font.Size = 10
$mid.Font.Bold = $true

$post = $c2.Characters($boldStart + $boldLen, 999)
$post.Font.Name = "Open Sans"
$post.Font.Size = 10
$post.Font.Bold = $false

# Propagate the exact same (shared) rich-text value to C3:C6 via copy/paste
# so they reuse one shared-string entry instead of independent duplicates.
foreach ($cellRef in @("C3", "C4", "C5", "C6")) {
    $c2.Copy()
    $ws.Range($cellRef).PasteSpecial()
}
$excel.CutCopyMode = $false

# --- 3. C7 ("Hände") loses the special "Open Sans" tool-column font and
#        reverts to the workbook default font, still centred vertically. ---
$c7 = $ws.Range("C7")
$c7.Font.Name = "Calibri"
$c7.Font.Size = 11

# --- 4. Ablauf row: normalise the placeholder text to "[Image_anguss]" ---
$ws.Range("B12").Value = "[Image_anguss]"
$ws.Range("D12").Value = "[Image_anguss]"

# --- 5. Row-height adjustments ---
$ws.Range("A1:D1").RowHeight = 13.8
$ws.Range("A2:D6").RowHeight = 16
$ws.Range("A7:D12").RowHeight = 13.8

# --- 6. Selection moves from D13 to A1:D12 (active cell A1) ---
$ws.Range("A1:D12").Select()
